$d = $word.ActiveDocument

function New-WordXmlPayload($innerBodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerBodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# --- Paragraph 1: "First checking wih Gibhub." -> add proofErr markers around
#     the misspelled / grammar-flagged words, keeping the same visible text.
$para1Body = '<w:p>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t xml:space="preserve">First checking </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>wih</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Gibhub</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '</w:p>'
$p1 = $d.Paragraphs.Item(1).Range
$p1.InsertXML((New-WordXmlPayload $para1Body))

# --- Paragraph 2: "ha ha ha" -> add proofErr markers, same visible text.
$para2Body = '<w:p>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>ha</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>ha</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>ha</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'
$p2 = $d.Paragraphs.Item(2).Range
$p2.InsertXML((New-WordXmlPayload $para2Body))

# --- Paragraph 4 ("Commit by Duc My Nguyen Second") gets a new sibling
#     paragraph after it containing the GitHub wiki link; the _GoBack
#     bookmark moves from paragraph 4 onto the new paragraph.
#     To reliably grow the body by one paragraph at the very end of the
#     document (where Word keeps an implicit, undeletable final mark),
#     insert two paragraph breaks first so our target paragraph is no
#     longer "last", fill both paragraphs via InsertXML, then remove the
#     now-redundant empty trailing paragraph that Word insists on keeping.
$p4 = $d.Paragraphs.Item(4).Range
$p4.InsertParagraphAfter()
$p5 = $d.Paragraphs.Item(5).Range
$p5.InsertParagraphAfter()

$para4Body = '<w:p><w:r><w:t>Commit by Duc My Nguyen Second</w:t></w:r></w:p>'
$p4 = $d.Paragraphs.Item(4).Range
$p4.InsertXML((New-WordXmlPayload $para4Body))

$para5Body = '<w:p><w:r><w:t>https://github.com/Kunena/Kunena-Forum/wiki/Create-a-new-branch-with-git-and-manage-branches</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$p5 = $d.Paragraphs.Item(5).Range
$p5.InsertXML((New-WordXmlPayload $para5Body))

# Remove the leftover empty paragraph Word kept at the very end.
$secondToLast = $d.Paragraphs.Item($d.Paragraphs.Count - 1).Range
$tailRange = $d.Range($secondToLast.End, $d.Content.End)
$tailRange.Delete()
